# feat: add 2022-Q4 data
#
# - Insert a new worksheet "2022-Q4" right after "总计" (i.e. as the 2nd
#   sheet, pushing the existing quarterly sheets down by one slot) holding
#   the new quarter's fund-holding detail rows.
# - Update the "总计" (summary) sheet with a new leading data row for
#   2022-Q4 and shift the previously-existing rows down by one.

$wb = $excel.ActiveWorkbook

function Set-TextCell($cell, $val) {
    # Force the cell to be written as text (Excel would otherwise happily
    # coerce strings like "012419" or "1.47" into numbers and lose the
    # leading zero / change type), then drop back to the default "Normal"
    # style so no stray number-format style is left behind on the cell.
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.Style = "Normal"
}

# ---------------------------------------------------------------------
# 1. Insert the new "2022-Q4" worksheet right after "总计".
#
# We duplicate the existing "2022-Q3" sheet (Worksheet.Copy) rather than
# adding a brand new blank sheet: this way the new sheet naturally inherits
# the exact same layout/column styles as every other quarterly detail
# sheet, and we only need to overwrite the data cells below the header.
# ---------------------------------------------------------------------
$wsTotal = $wb.Worksheets.Item(1)
$wsQ3Old = $wb.Worksheets.Item(2)   # current "2022-Q3" sheet

$wsQ3Old.Copy($wsQ3Old)             # places the duplicate right before it
$wsQ4 = $wb.Worksheets.Item(2)
$wsQ4.Name = "2022-Q4"

$q4rows = @(
    @(0, "012419", "天弘国证建材指数C", "0.48", "94.64", "3.20", "0.0154", 9),
    @(1, "012405", "天弘国证建材指数A", "0.14", "94.64", "3.20", "0.0045", 9)
)

$r = 2
foreach ($row in $q4rows) {
    $wsQ4.Cells.Item($r, 1).Value = $row[0]
    Set-TextCell $wsQ4.Cells.Item($r, 2) $row[1]
    Set-TextCell $wsQ4.Cells.Item($r, 3) $row[2]
    Set-TextCell $wsQ4.Cells.Item($r, 4) $row[3]
    Set-TextCell $wsQ4.Cells.Item($r, 5) $row[4]
    Set-TextCell $wsQ4.Cells.Item($r, 6) $row[5]
    Set-TextCell $wsQ4.Cells.Item($r, 7) $row[6]
    $wsQ4.Cells.Item($r, 8).Value = $row[7]
    $r++
}

# ---------------------------------------------------------------------
# 2. Update the "总计" (summary) sheet: add the 2022-Q4 row at the top of
#    the data and shift everything else down by one row.
# ---------------------------------------------------------------------
$summaryRows = @(
    @(0, "2022-Q4", 2, 0.02),
    @(1, "2022-Q3", 2, 0.04),
    @(2, "2022-Q1", 3, 0.02),
    @(3, "2021-Q4", 8, 1.29),
    @(4, "2021-Q3", 5, 1.21),
    @(5, "2021-Q2", 2, 0.11),
    @(6, "2021-Q1", 3, 0.53),
    @(7, "2020-Q4", 5, 3.15)
)

$r = 2
foreach ($row in $summaryRows) {
    $wsTotal.Cells.Item($r, 1).Value = $row[0]
    Set-TextCell $wsTotal.Cells.Item($r, 2) $row[1]
    $wsTotal.Cells.Item($r, 3).Value = $row[2]
    $wsTotal.Cells.Item($r, 4).Value = $row[3]
    $r++
}

# Row 9 is brand new on this sheet -- copy A2's format (the bold/bordered
# id-column look) onto the new A9 first, then overwrite its value.
$wsTotal.Cells.Item(2, 1).Copy()
$wsTotal.Cells.Item(9, 1).PasteSpecial(-4122)
$wsTotal.Cells.Item(9, 1).Value = 7

# Restore the original active tab ("2020-Q4", the last sheet) and a clean
# A1 selection on every sheet, matching the saved view state.
foreach ($s in $wb.Worksheets) {
    $s.Range("A1").Select() | Out-Null
}
$wb.Worksheets.Item($wb.Worksheets.Count).Activate()
